$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd() -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Change 1: remove the stray _GoBack bookmark left on the empty paragraph
# right after the document title.
# ---------------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    if ($goBack -ne $null) {
        $goBack.Delete()
    }
} catch {
    # no _GoBack bookmark present - nothing to remove
}

# ---------------------------------------------------------------------------
# Change 2: add a new concluding paragraph right after the
# "... invest in the technological market." paragraph.
# ---------------------------------------------------------------------------
$anchorPara = Find-ParagraphByText $d "We will take into account as many factors as possible for our simulations, with the purpose of obtaining the most credible results and discover how to best invest in the technological market."

$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$insPoint = $newPara.Range
$insPoint.Collapse(1)
$paraStart = $insPoint.Start

$leadIn = "In short, we propose a multiagent system for investment simulation in the technological market"
$fullText = $leadIn + " where our agents are the companies, the company managers and Investors."
$insPoint.InsertAfter($fullText)

# Make the whole paragraph (including its paragraph mark) bold, then turn
# bold back off for the spans that must stay regular weight (the bold runs
# are "the companies", "the company managers " and "Investors.").
$newPara.Range.Font.Bold = $true

$boldStart1 = $fullText.IndexOf("the companies")
$boldEnd1 = $boldStart1 + "the companies".Length
$boldStart2 = $fullText.IndexOf("the company managers ")
$boldEnd2 = $boldStart2 + "the company managers ".Length
$boldStart3 = $fullText.IndexOf("Investors.")
$boldEnd3 = $boldStart3 + "Investors.".Length

$d.Range($paraStart, $paraStart + $boldStart1).Font.Bold = $false
$d.Range($paraStart + $boldEnd1, $paraStart + $boldStart2).Font.Bold = $false
$d.Range($paraStart + $boldEnd2, $paraStart + $boldStart3).Font.Bold = $false

# Re-insert the _GoBack bookmark right after "... technological market"
# (collapsed range, i.e. zero-length) inside the new paragraph.
$bmPos = $paraStart + $leadIn.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# Change 3: merge the "Number of companies" bullet into the
# "Number of investors" bullet, turning it into
# "Number of investors and companies".
# ---------------------------------------------------------------------------
$investorsPara = Find-ParagraphByText $d "Number of investors"
$mergeStart = $investorsPara.Range.End
$paraMarkRange = $d.Range($mergeStart - 1, $mergeStart)
$paraMarkRange.Delete()

$investorsPara = Find-ParagraphByText $d "Number of investorsNumber of companies"
$start = $investorsPara.Range.Start
$companiesRange = $d.Range($start + 19, $start + 19 + 19)
$companiesRange.Text = " and companies"
